# Apply numeric corrections to the Leve profit-tracking sheets.
# Generated from the authoritative cell-level diff (values only; no formulas involved).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 33
$ws.Range("H33").Value = 350.85184
$ws.Range("I33").Value = 357.95834
$ws.Range("J33").Value = 294
$ws.Range("K33").Value = 357.95834
$ws.Range("L33").Value = 294
$ws.Range("M33").Value = -128.95834
$ws.Range("N33").Value = -752
# row 68
$ws.Range("H68").Value = 30000
$ws.Range("J68").Value = 30000
$ws.Range("L68").Value = 30000
$ws.Range("N68").Value = -31498
# row 71
$ws.Range("H71").Value = 30000
$ws.Range("J71").Value = 30000
$ws.Range("L71").Value = 90000
$ws.Range("N71").Value = -97488
# row 107
$ws.Range("H107").Value = 7873.077
$ws.Range("I107").Value = 12618.125
$ws.Range("J107").Value = 281
$ws.Range("K107").Value = 12618.125
$ws.Range("L107").Value = 281
$ws.Range("M107").Value = -10698.125
$ws.Range("N107").Value = -4121
# row 121
$ws.Range("H121").Value = 821.43634
$ws.Range("J121").Value = 814.7925
$ws.Range("L121").Value = 2444.3775
$ws.Range("N121").Value = -5938.377500000001
# row 135
$ws.Range("H135").Value = 1203.8
$ws.Range("J135").Value = 3000
$ws.Range("L135").Value = 27000
$ws.Range("N135").Value = -32070
# row 141
$ws.Range("H141").Value = 1763.2354
$ws.Range("I141").Value = 1685.9375
$ws.Range("K141").Value = 5057.8125
$ws.Range("M141").Value = 122.1875

$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 13940.779
$ws.Range("I32").Value = 13480.808
$ws.Range("K32").Value = 13480.808
$ws.Range("M32").Value = -13193.808
# row 45
$ws.Range("H45").Value = 1684727.8
$ws.Range("I45").Value = 1895168.8
$ws.Range("J45").Value = 1200
$ws.Range("K45").Value = 1895168.8
$ws.Range("L45").Value = 1200
$ws.Range("M45").Value = -1894791.8
$ws.Range("N45").Value = -1954
# row 122
$ws.Range("H122").Value = 7170.4
$ws.Range("I122").Value = 7489.6523
$ws.Range("J122").Value = 3499
$ws.Range("K122").Value = 22468.9569
$ws.Range("L122").Value = 10497
$ws.Range("M122").Value = -20018.9569
$ws.Range("N122").Value = -15397
# row 132
$ws.Range("H132").Value = 1872.4386
$ws.Range("I132").Value = 1392.9286
$ws.Range("J132").Value = 3215.0667
$ws.Range("K132").Value = 4178.7858
$ws.Range("L132").Value = 9645.2001
$ws.Range("M132").Value = -1648.7858
$ws.Range("N132").Value = -14705.2001

$ws = $wb.Worksheets.Item("BSM")
# row 134
$ws.Range("H134").Value = 3274.7969
$ws.Range("I134").Value = 2226.325
$ws.Range("J134").Value = 5022.25
$ws.Range("K134").Value = 6678.974999999999
$ws.Range("L134").Value = 15066.75
$ws.Range("M134").Value = -4143.974999999999
$ws.Range("N134").Value = -20136.75

$ws = $wb.Worksheets.Item("CRP")
# row 94
$ws.Range("H94").Value = 3319.4666
$ws.Range("I94").Value = 2129.4614
$ws.Range("J94").Value = 4229.4707
$ws.Range("K94").Value = 2129.4614
$ws.Range("L94").Value = 4229.4707
$ws.Range("M94").Value = -1678.4614
$ws.Range("N94").Value = -5131.4707
# row 132
$ws.Range("H132").Value = 2127.1765
$ws.Range("I132").Value = 1691.6154
$ws.Range("J132").Value = 3542.75
$ws.Range("K132").Value = 5074.8462
$ws.Range("L132").Value = 10628.25
$ws.Range("M132").Value = -2544.8462
$ws.Range("N132").Value = -15688.25

$ws = $wb.Worksheets.Item("CUL")
# row 23
$ws.Range("H23").Value = 104.083336
$ws.Range("I23").Value = 32.25
$ws.Range("J23").Value = 140
$ws.Range("K23").Value = 96.75
$ws.Range("L23").Value = 420
$ws.Range("M23").Value = 138.25
$ws.Range("N23").Value = -890
# row 75
$ws.Range("M75").Value = -472
$ws.Range("H75").Value = 3789.2
$ws.Range("I75").Value = 490
$ws.Range("J75").Value = 4155.778
$ws.Range("K75").Value = 1470
$ws.Range("L75").Value = 12467.334
$ws.Range("N75").Value = -14463.334
# row 78
$ws.Range("M78").Value = 582
$ws.Range("H78").Value = 3789.2
$ws.Range("I78").Value = 490
$ws.Range("J78").Value = 4155.778
$ws.Range("K78").Value = 4410
$ws.Range("L78").Value = 37402.002
$ws.Range("N78").Value = -47386.002
# row 81
$ws.Range("H81").Value = 5904.6665
$ws.Range("J81").Value = 6910.6
$ws.Range("L81").Value = 20731.8
$ws.Range("N81").Value = -22977.8
# row 84
$ws.Range("H84").Value = 5904.6665
$ws.Range("J84").Value = 6910.6
$ws.Range("L84").Value = 62195.4
$ws.Range("N84").Value = -73427.39999999999
# row 134
$ws.Range("H134").Value = 5030.6895
$ws.Range("I134").Value = 2558.2354
$ws.Range("J134").Value = 8533.333000000001
$ws.Range("K134").Value = 7674.706200000001
$ws.Range("L134").Value = 25599.999
$ws.Range("M134").Value = -2604.706200000001
$ws.Range("N134").Value = -35739.999
# row 139
$ws.Range("H139").Value = 1929.6774
$ws.Range("I139").Value = 1096.6666
$ws.Range("J139").Value = 4785.7144
$ws.Range("K139").Value = 3289.9998
$ws.Range("L139").Value = 14357.1432
$ws.Range("M139").Value = 1850.0002
$ws.Range("N139").Value = -24637.1432
# row 140
$ws.Range("H140").Value = 2561.577
$ws.Range("I140").Value = 1440.0667
$ws.Range("J140").Value = 4090.9092
$ws.Range("K140").Value = 4320.2001
$ws.Range("L140").Value = 12272.7276
$ws.Range("M140").Value = 859.7999
$ws.Range("N140").Value = -22632.7276

$ws = $wb.Worksheets.Item("GSM")
# row 57
$ws.Range("H57").Value = 10518.091
$ws.Range("J57").Value = 15550
$ws.Range("L57").Value = 15550
$ws.Range("N57").Value = -17190
# row 122
$ws.Range("H122").Value = 2662.818
$ws.Range("I122").Value = 2131.077
$ws.Range("J122").Value = 3430.889
$ws.Range("K122").Value = 6393.231000000001
$ws.Range("L122").Value = 10292.667
$ws.Range("M122").Value = -3943.231000000001
$ws.Range("N122").Value = -15192.667
# row 132
$ws.Range("H132").Value = 4652.222
$ws.Range("I132").Value = 5165.1177
$ws.Range("J132").Value = 3066.9092
$ws.Range("K132").Value = 15495.3531
$ws.Range("L132").Value = 9200.7276
$ws.Range("M132").Value = -12965.3531
$ws.Range("N132").Value = -14260.7276

$ws = $wb.Worksheets.Item("LTW")
# row 4
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
# row 22
$ws.Range("H22").Value = 1295.3
$ws.Range("I22").Value = 881.6667
$ws.Range("J22").Value = 1472.5714
$ws.Range("K22").Value = 881.6667
$ws.Range("L22").Value = 1472.5714
$ws.Range("M22").Value = -586.6667
$ws.Range("N22").Value = -2062.5714
# row 27
$ws.Range("H27").Value = 1295.3
$ws.Range("I27").Value = 881.6667
$ws.Range("J27").Value = 1472.5714
$ws.Range("K27").Value = 881.6667
$ws.Range("L27").Value = 1472.5714
$ws.Range("M27").Value = -774.6667
$ws.Range("N27").Value = -1686.5714
# row 28
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
# row 37
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
# row 40
$ws.Range("H40").Value = 4712.2354
$ws.Range("I40").Value = 7663.5
$ws.Range("J40").Value = 2088.889
$ws.Range("K40").Value = 7663.5
$ws.Range("L40").Value = 2088.889
$ws.Range("M40").Value = -7527.5
$ws.Range("N40").Value = -2360.889
# row 122
$ws.Range("H122").Value = 13033.733
$ws.Range("I122").Value = 13011.6
$ws.Range("J122").Value = 13044.8
$ws.Range("K122").Value = 39034.8
$ws.Range("L122").Value = 39134.39999999999
$ws.Range("M122").Value = -36584.8
$ws.Range("N122").Value = -44034.39999999999
# row 132
$ws.Range("H132").Value = 8338809.5
$ws.Range("I132").Value = 4298.3784
$ws.Range("J132").Value = 21746500
$ws.Range("K132").Value = 12895.1352
$ws.Range("L132").Value = 65239500
$ws.Range("M132").Value = -10365.1352
$ws.Range("N132").Value = -65244560
# row 136
$ws.Range("H136").Value = 6148.8276
$ws.Range("I136").Value = 2629.4167
$ws.Range("K136").Value = 7888.250100000001
$ws.Range("M136").Value = -5338.250100000001

$ws = $wb.Worksheets.Item("WVR")
# row 21
$ws.Range("M21").Value = -8765
$ws.Range("H21").Value = 11166.8
$ws.Range("I21").Value = 9000
$ws.Range("J21").Value = 11708.5
$ws.Range("K21").Value = 9000
$ws.Range("L21").Value = 11708.5
$ws.Range("N21").Value = -12178.5
# row 35
$ws.Range("M35").Value = -8710
$ws.Range("H35").Value = 11166.8
$ws.Range("I35").Value = 9000
$ws.Range("J35").Value = 11708.5
$ws.Range("K35").Value = 9000
$ws.Range("L35").Value = 11708.5
$ws.Range("N35").Value = -12288.5
# row 62
$ws.Range("H62").Value = 5367.2173
$ws.Range("I62").Value = 5269.8
$ws.Range("J62").Value = 5442.154
$ws.Range("K62").Value = 5269.8
$ws.Range("L62").Value = 5442.154
$ws.Range("M62").Value = -4645.8
$ws.Range("N62").Value = -6690.154
# row 65
$ws.Range("H65").Value = 5367.2173
$ws.Range("I65").Value = 5269.8
$ws.Range("J65").Value = 5442.154
$ws.Range("K65").Value = 26349
$ws.Range("L65").Value = 27210.77
$ws.Range("M65").Value = -23229
$ws.Range("N65").Value = -33450.77
# row 107
$ws.Range("H107").Value = 2290.4
$ws.Range("I107").Value = 10002
$ws.Range("J107").Value = 362.5
$ws.Range("K107").Value = 30006
$ws.Range("L107").Value = 1087.5
$ws.Range("M107").Value = -28086
$ws.Range("N107").Value = -4927.5
# row 132
$ws.Range("H132").Value = 2633.0303
$ws.Range("I132").Value = 2073.926
$ws.Range("J132").Value = 5149
$ws.Range("K132").Value = 6221.778
$ws.Range("L132").Value = 15447
$ws.Range("M132").Value = -3691.778
$ws.Range("N132").Value = -20507
# row 136
$ws.Range("H136").Value = 1761.8223
$ws.Range("I136").Value = 935.9697
$ws.Range("K136").Value = 2807.9091
$ws.Range("M136").Value = -257.9090999999999
